# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds the "K" values (header in G1 is "K"). Update the per-row K
# values for the 8 data rows (rows 2-9) to the newly regenerated numbers.
$newK = @{
    2 = 4
    3 = 1
    4 = 0
    5 = 1
    6 = 0
    7 = 1
    8 = 1
    9 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
